# Update the cryptos price/volume table with refreshed values scraped by the
# GitHub Actions workflow. Values are written as text (matching the original
# inlineStr cell type) so that numeric-looking strings such as "238.62" are
# not silently reinterpreted by Excel as numbers, and percentage strings keep
# their exact surrounding whitespace.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param(
        $ws,
        [string]$addr,
        [string]$val
    )
    $cell = $ws.Range($addr)
    # Prefixing with an apostrophe forces Excel to store the value as text
    # even when it looks like a number (e.g. "238.62" or "0.100").
    $cell.Value = "'" + $val
    # Re-apply the default "Normal" style so no stray number-format / quote
    # prefix style is left behind on the cell.
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" '43.989.08'
Set-CellText $ws "E2" '  -1.14%  '
Set-CellText $ws "D3" '2.349.63'
Set-CellText $ws "E3" '  -0.81%  '
Set-CellText $ws "E4" '  +0.10%  '
Set-CellText $ws "D5" '0.679'
Set-CellText $ws "E5" '  +0.12%  '
Set-CellText $ws "D6" '238.62'
Set-CellText $ws "E6" '  -0.19%  '
Set-CellText $ws "D7" '73.72'
Set-CellText $ws "E7" '  +0.74%  '
Set-CellText $ws "E8" '  +0.06%  '
Set-CellText $ws "D9" '0.593'
Set-CellText $ws "E9" '  +8.86%  '
Set-CellText $ws "D10" '0.100'
Set-CellText $ws "E10" '  -2.88%  '
Set-CellText $ws "E11" '  -0.35%  '
Set-CellText $ws "D12" '32.38'
Set-CellText $ws "E12" '  +10.24%  '
Set-CellText $ws "D13" '7.28'
Set-CellText $ws "E13" '  +8.56%  '
Set-CellText $ws "E14" '  +0.30%  '
Set-CellText $ws "D15" '2.697.70'
Set-CellText $ws "E15" '  -0.66%  '
Set-CellText $ws "D16" '16.55'
Set-CellText $ws "E16" '  -1.67%  '
Set-CellText $ws "D17" '0.897'
Set-CellText $ws "E17" '  -0.81%  '
Set-CellText $ws "D18" '2.353.20'
Set-CellText $ws "E18" '  -0.40%  '
Set-CellText $ws "D19" '43.847.55'
Set-CellText $ws "E19" '  -1.22%  '
Set-CellText $ws "E20" '  -2.84%  '
Set-CellText $ws "D21" '6.73'
Set-CellText $ws "E21" '  +4.35%  '
Set-CellText $ws "D22" '76.81'
Set-CellText $ws "E22" '  -1.41%  '
Set-CellText $ws "D23" '258.99'
Set-CellText $ws "E23" '  +1.34%  '
Set-CellText $ws "D24" '1.95'
Set-CellText $ws "E24" '  +22.86%  '
Set-CellText $ws "E25" '  +0.12%  '
Set-CellText $ws "D26" '3.66'
Set-CellText $ws "E26" '  -2.49%  '
Set-CellText $ws "E27" '  -1.90%  '
Set-CellText $ws "D28" '10.71'
Set-CellText $ws "E28" '  +2.22%  '
Set-CellText $ws "E29" '  -0.51%  '
Set-CellText $ws "D30" '22.62'
Set-CellText $ws "E30" '  +0.41%  '
Set-CellText $ws "D31" '175.54'
Set-CellText $ws "E31" '  +1.33%  '
Set-CellText $ws "E32" '  -3.54%  '
Set-CellText $ws "E33" '  +2.24%  '
Set-CellText $ws "D34" '0.0758'
Set-CellText $ws "E34" '  +2.36%  '
Set-CellText $ws "D35" '5.20'
Set-CellText $ws "E35" '  -0.04%  '
Set-CellText $ws "D36" '5.47'
Set-CellText $ws "E36" '  +5.07%  '
Set-CellText $ws "D37" '3.74'
Set-CellText $ws "E37" '  -4.78%  '
Set-CellText $ws "E38" '  -3.81%  '
Set-CellText $ws "D39" '6.28'
Set-CellText $ws "D40" '0.0277'
Set-CellText $ws "E40" '  +1.76%  '
Set-CellText $ws "E41" '  +12.28%  '
Set-CellText $ws "D42" '0.205'
Set-CellText $ws "E42" '  +11.84%  '
Set-CellText $ws "D43" '18.89'
Set-CellText $ws "E43" '  -4.40%  '
Set-CellText $ws "E44" '  -0.07%  '
Set-CellText $ws "D45" '8.94'
Set-CellText $ws "E45" '  +0.67%  '
Set-CellText $ws "E46" '  +4.51%  '
Set-CellText $ws "D47" '2.51'
Set-CellText $ws "E47" '  +6.48%  '
Set-CellText $ws "D48" '57.74'
Set-CellText $ws "E48" '  +9.55%  '
Set-CellText $ws "E49" '  -1.73%  '
Set-CellText $ws "D50" '1.17'
Set-CellText $ws "E50" '  -0.16%  '
Set-CellText $ws "D51" '99.75'
Set-CellText $ws "E51" '  +1.06%  '
